# LOB1024.xlsx update
# - "Objetivos:" value replaced with the docente string
# - Old "Docentes responsáveis:" value row (row 13, blank label) is removed entirely
#   (its content was promoted into the "Objetivos:" row above), which shifts every
#   subsequent row up by one.
# - "Programa resumido:" value replaced with "Semestral"
# - "Programa:" value replaced with "01/01/2018"
# - "Método:" value replaced with the docente string
# - "Critério:" / "Norma de recuperação:" / "Bibliografia:" values each take on the
#   value that used to belong to the row below them; the old Bibliografia text block
#   is dropped entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "Objetivos:" (row 10) gets the docente string instead of the old objective text.
$ws.Range("B10").Value = "5840650 - Janaína Ferreira Batista"
$ws.Range("C10").Value = "5840650 - Janaína Ferreira Batista"

# 2. Remove the stray row 13 (label-less row holding the old docente string) -
#    everything below shifts up by one row.
$ws.Rows(13).Delete()

# 3. "Programa resumido:" (now row 13) becomes "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# 4. "Programa:" (now row 15) becomes "01/01/2018".
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

# 5. "Método:" (now row 18) gets the docente string.
$ws.Range("B18").Value = "5840650 - Janaína Ferreira Batista"
$ws.Range("C18").Value = "5840650 - Janaína Ferreira Batista"

# 6. "Critério:" (now row 19) takes the old "Método:" evaluation text.
$ws.Range("B19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."
$ws.Range("C19").Value = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# 7. "Norma de recuperação:" (now row 20) takes the old "Critério:" text.
$ws.Range("B20").Value = "NF≥ 5,0."
$ws.Range("C20").Value = "NF≥ 5,0."

# 8. "Bibliografia:" (now row 21) takes the old "Norma de recuperação:" text.
$ws.Range("B21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
$ws.Range("C21").Value = "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada."
